$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Nombre Problema",
    "Alpha=0",
    "Punto original (x_1,x_2,..., x_n,y_1,y_2,...,y_m)",
    "Valor Objetivo Nivel Superior",
    "Punto obtenido ahora julia (x_1,x_2,..., x_n,y_1,y_2,...,y_m)",
    "Valor Objetivo Nivel Superior Obtenido por Julia",
    "Optimizador"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws.Columns.Item(1).ColumnWidth = 20.0
$ws.Columns.Item(2).ColumnWidth = 27.5
$ws.Columns.Item(3).ColumnWidth = 43.5
$ws.Columns.Item(4).ColumnWidth = 26.6666666666667
$ws.Columns.Item(5).ColumnWidth = 54.8333333333333
$ws.Columns.Item(6).ColumnWidth = 44.3333333333333
$ws.Columns.Item(7).ColumnWidth = 11.6666666666667

$ws.Rows("3:10").Select()
